# Update cryptos list values per data refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cell, $value) {
    $rng = $ws.Range($cell)
    $rng.NumberFormat = "@"
    $rng.Value = $value
}

Set-TextCell "D2" "43.030.72"
Set-TextCell "E2" "  -0.35%  "

Set-TextCell "D3" "2.303.93"
Set-TextCell "E3" "  -0.74%  "

Set-TextCell "E4" "  -0.02%  "

Set-TextCell "D5" "299.87"

Set-TextCell "D6" "98.47"
Set-TextCell "E6" "  -1.54%  "

Set-TextCell "D7" "0.518"
Set-TextCell "E7" "  +2.17%  "

Set-TextCell "E8" "  -0.05%  "

Set-TextCell "E9" "  -1.50%  "

Set-TextCell "D10" "36.45"
Set-TextCell "E10" "  -0.53%  "

Set-TextCell "D11" "0.0789"
Set-TextCell "E11" "  -0.36%  "

Set-TextCell "B12" "Chainlink"
Set-TextCell "C12" "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextCell "D12" "17.91"
Set-TextCell "E12" "  +0.01%  "

Set-TextCell "B13" "TRON"
Set-TextCell "C13" "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
Set-TextCell "D13" "0.117"
Set-TextCell "E13" "  +0.59%  "

Set-TextCell "E14" "  -2.21%  "

Set-TextCell "D15" "2.662.25"
Set-TextCell "E15" "  -0.86%  "

Set-TextCell "D16" "2.304.32"
Set-TextCell "E16" "  -1.90%  "

Set-TextCell "E17" "  -1.87%  "

Set-TextCell "D18" "42.999.27"
Set-TextCell "E18" "  -0.31%  "

Set-TextCell "D19" "12.63"

Set-TextCell "E20" "  +0.30%  "

Set-TextCell "D21" "6.13"
Set-TextCell "E21" "  -1.03%  "

Set-TextCell "D22" "68.45"
Set-TextCell "E22" "  +0.39%  "

Set-TextCell "D23" "242.14"
Set-TextCell "E23" "  +2.31%  "

Set-TextCell "D24" "2.15"
Set-TextCell "E24" "  -1.19%  "

Set-TextCell "E25" "  +0.08%  "

Set-TextCell "D26" "2.44"
Set-TextCell "E26" "  -1.18%  "

Set-TextCell "E27" "  -0.18%  "

Set-TextCell "D28" "25.13"
Set-TextCell "E28" "  -0.92%  "

Set-TextCell "D29" "166.45"
Set-TextCell "E29" "  -1.31%  "

Set-TextCell "D30" "2.04"
Set-TextCell "E30" "  -0.36%  "

Set-TextCell "D31" "9.10"
Set-TextCell "E31" "  -0.86%  "

Set-TextCell "D32" "33.26"
Set-TextCell "E32" "  -3.85%  "

Set-TextCell "E33" "  +0.02%  "

Set-TextCell "D34" "5.03"
Set-TextCell "E34" "  -2.99%  "

Set-TextCell "E35" "  +0.78%  "

Set-TextCell "D36" "17.77"
Set-TextCell "E36" "  +1.05%  "

Set-TextCell "E37" "  -0.26%  "

Set-TextCell "E38" "  -0.87%  "

Set-TextCell "E39" "  -1.10%  "

Set-TextCell "E40" "  -1.75%  "

Set-TextCell "D41" "2.78"
Set-TextCell "E41" "  +0.12%  "

Set-TextCell "E42" "  +0.58%  "

Set-TextCell "D43" "2.000.99"
Set-TextCell "E43" "  +0.20%  "

Set-TextCell "E44" "  -1.85%  "

Set-TextCell "D45" "2.19"
Set-TextCell "E45" "  -3.40%  "

Set-TextCell "D46" "10.22"
Set-TextCell "E46" "  +0.92%  "

Set-TextCell "D47" "17.45"
Set-TextCell "E47" "  -2.39%  "

Set-TextCell "E48" "  -3.34%  "

Set-TextCell "D49" "53.85"
Set-TextCell "E49" "  -2.49%  "

Set-TextCell "D50" "2.528.30"
Set-TextCell "E50" "  -0.85%  "

Set-TextCell "B51" "BitcoinSV"
Set-TextCell "C51" "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
Set-TextCell "D51" "72.71"
Set-TextCell "E51" "  +1.20%  "
